$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.434268236160278
$ws.Range("B1").Value = 2.251269578933716
$ws.Range("C1").Value = 1.693834662437439
$ws.Range("D1").Value = 1.805990934371948
$ws.Range("E1").Value = 1.605077266693115
